$d = $word.ActiveDocument

# --- Step 1: delete the trailing image-only paragraphs (23.jpg, 22.jpg, 21.jpg, 20.jpg) ---
# These are paragraphs 5-8 in the original document; delete from the end backwards
# so earlier indices stay valid.
$d.Paragraphs.Item(8).Range.Delete()
$d.Paragraphs.Item(7).Range.Delete()
$d.Paragraphs.Item(6).Range.Delete()
$d.Paragraphs.Item(5).Range.Delete()

# --- Step 2: replace the 24.jpg image-only paragraph (now paragraph 4) with text "(1) |" ---
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Text = "(1) |"

# --- Step 3: fix typo "between X and the length" -> "between X and! the length" ---
$found = $d.Content.Find.Execute("between X and the length", $true, $false, $false, $false, $false,
                         $true, 1, $false, "between X and! the length", 2)

# --- Step 4: remove the trailing line-break + "|" at the end of the "Based on..." paragraph ---
$p3 = $d.Paragraphs.Item(3)
$pStart = $p3.Range.Start
$fullText = $p3.Range.Text
$idx = $fullText.IndexOf("screen?")
$cutStart = $pStart + $idx + 7      # position right after "screen?"
$cutEnd = $p3.Range.End - 1         # exclude the paragraph mark itself
$trailRange = $d.Range($cutStart, $cutEnd)
$trailRange.Delete()

# --- Step 5: delete the 25.jpg image-only paragraph (paragraph 2); it merges into paragraph 1 ---
$d.Paragraphs.Item(2).Range.Delete()

# --- Step 6: append a line break + "|" to the end of paragraph 1 ("18. Study the set-up below.") ---
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertAfter([char]11 + "|")
